$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Team 10 Mechanik"

# New order rows (10-13)
$ws.Cells.Item(10,1).Value = "distrelec.ch"
$ws.Cells.Item(10,2).Value = 1
$ws.Cells.Item(10,3).Value = "182-49-575"
$ws.Cells.Item(10,4).Value = "Lötpaste"
$ws.Cells.Item(10,5).Value = 17.4
$ws.Range("F10").Formula = "=B10*E10"

$ws.Cells.Item(11,1).Value = "distrelec.ch"
$ws.Cells.Item(11,2).Value = 1
$ws.Cells.Item(11,3).Value = "182-49-574"
$ws.Cells.Item(11,4).Value = "Dosiernadel"
$ws.Cells.Item(11,5).Value = 1.08
$ws.Range("F11").Formula = "=B11*E11"

$ws.Cells.Item(12,1).Value = "distrelec.ch"
$ws.Cells.Item(12,2).Value = 1
$ws.Cells.Item(12,3).Value = "300-72-943"
$ws.Cells.Item(12,4).Value = "Schlitten für lineare Führungsschienen 9 mm, Vorgeschmiert ja, MNN 9-G3-LS, Schneeberger"
$ws.Cells.Item(12,5).Value = 29.67
$ws.Range("F12").Formula = "=B12*E12"

$ws.Cells.Item(13,1).Value = "distrelec.ch"
$ws.Cells.Item(13,2).Value = 1
$ws.Cells.Item(13,3).Value = "300-72-923"
$ws.Cells.Item(13,4).Value = "Lineare Führungsschiene 9/275 mm, MN 9-275-G3-V0, Schneeberger"
$ws.Cells.Item(13,5).Value = 45.05
$ws.Range("F13").Formula = "=B13*E13"

# Clear old placeholder formulas in F14:F30 (they become blank)
$ws.Range("F14:F30").ClearContents()

# Column C width
$ws.Columns.Item(3).ColumnWidth = 21.33203125

# Selection
$ws.Range("D17").Select()
